# BaoCao_SystemTest_Login.xlsx edit
# - Update the "sai mat khau" (ST_LOGIN_03) row wording / messages
# - Insert a brand-new "ST_LOGIN_04" (tai khoan khong ton tai) row after it
# - Append "chuyen huong" wording to the two existing success rows
# - Widen columns D and E to fit the longer text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Edit the existing ST_LOGIN_03 row (row 2) in place ---------------
$ws.Range("C2").Value = "1. Nhập user/sai_pass`n2. Click Login`n3. Check URL & Thông báo lỗi"
$ws.Range("D2").Value = "User: admin, Pass: sai_pass_roi_nhe"
$ws.Range("E2").Value = "Ở lại trang Login & Hiện lỗi 'Mật khẩu không chính xác'"
$ws.Range("F2").Value = "URL: Login.jsp | Lỗi: Mật khẩu không chính xác!"

# --- 2. Insert a new row before row 3; it inherits formatting from above -
$ws.Rows("3:3").Insert()

# Fill the freshly inserted row 3 with the new ST_LOGIN_04 test case
$ws.Range("A3").Value = "ST_LOGIN_04"
$ws.Range("B3").Value = "Đăng nhập tài khoản không tồn tại"
$ws.Range("C3").Value = "1. Nhập user rác`n2. Nhập pass bất kỳ`n3. Click Login`n4. Check lỗi 'không tồn tại'"
$ws.Range("D3").Value = "User: user_ao_ma_canada"
$ws.Range("E3").Value = "Hiện lỗi 'Tài khoản không tồn tại'"
$ws.Range("F3").Value = "URL: Login.jsp | Lỗi: Tài khoản không tồn tại!"
$ws.Range("G3").Value = "PASS"

# --- 3. The old row 3 (ST_LOGIN_01) is now row 4; append wording ---------
$ws.Range("C4").Value = "1. Nhập admin/admin123`n2. Click Login`n3. Check URL chuyển hướng"

# --- 4. The old row 4 (ST_LOGIN_02) is now row 5; append wording ---------
$ws.Range("C5").Value = "1. Nhập user/user123`n2. Click Login`n3. Check URL chuyển hướng"

# --- 5. Re-fit the row heights so the embedded line breaks above don't ---
#        leave a stray explicit/custom row height behind.
$ws.Rows("1:5").AutoFit()

# --- 6. Widen columns D and E to fit the new/updated content -------------
# (COM ColumnWidth snaps to whole on-screen pixels -- 32.5 / 49.0 are the
#  closest achievable values to the saved XML widths of 33.40234375 /
#  49.7890625 characters.)
$ws.Columns("D").ColumnWidth = 32.5
$ws.Columns("E").ColumnWidth = 49.0
